$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.105.37"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.69%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.566.28"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.96%  "
$ws.Range("E4").Value = "  +0.61%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.491"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.48%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "21.92"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.45%  "
$ws.Range("E9").Value = "  +0.37%  "
$ws.Range("E10").Value = "  +0.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0863"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.82%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.788.59"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.91%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.523.57"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.77"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.68%  "
$ws.Range("E15").Value = "  -0.09%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.132.13"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.85%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.95"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.58%  "
$ws.Range("E18").Value = "  -0.63%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "214.75"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.94%  "
$ws.Range("E21").Value = "  +0.68%  "
$ws.Range("E22").Value = "  +1.32%  "
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.07"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.37%  "
$ws.Range("E26").Value = "  -0.51%  "
$ws.Range("E27").Value = "  +0.44%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.106"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.22%  "
$ws.Range("E29").Value = "  +0.61%  "
$ws.Range("E30").Value = "  +5.08%  "
$ws.Range("E31").Value = "  +0.52%  "
$ws.Range("E32").Value = "  +0.71%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.427.90"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.68%  "
$ws.Range("E35").Value = "  -4.88%  "
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("E37").Value = "  +2.21%  "
$ws.Range("E38").Value = "  +1.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.531"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.82"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.54%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.806"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.11%  "
$ws.Range("E42").Value = "  +0.73%  "
$ws.Range("E43").Value = "  +2.35%  "
$ws.Range("E44").Value = "  +0.62%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.52"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.13%  "
$ws.Range("E46").Value = "  +0.81%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.707.96"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.30%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.01"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.32%  "
$ws.Range("E49").Value = "  +1.63%  "
$ws.Range("E50").Value = "  -0.53%  "
$ws.Range("E51").Value = "  +0.26%  "
